$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3:22 down to 4:23
$ws.Rows(3).Insert()

# Make sure the label cell keeps the same style as the other label cells in column A
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Fill in the new row 3 with the new data point
$ws.Range("A3").Value = "2020-05-15 00:00:00_diff"
$ws.Range("B3").Value = -9.678192681000001
$ws.Range("C3").Value = 16.59095765
$ws.Range("D3").Value = -4.680304
$ws.Range("E3").Value = 1.859367
$ws.Range("F3").Value = 2.906937868
$ws.Range("G3").Value = -1.133092
$ws.Range("H3").Value = -1.294266
